$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9009016752243042
$ws.Range("B1").Value = 3.551664352416992
$ws.Range("C1").Value = 4.238682746887207
$ws.Range("D1").Value = 2.667722940444946
$ws.Range("E1").Value = 1.055494070053101
